$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "21.032.20"
$ws.Range("E2").Value = "  +3.14%  "

# Row 3
$ws.Range("D3").Value = "1.537.22"
$ws.Range("E3").Value = "  +5.06%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.51%  "

# Row 5
$ws.Range("D5").Value = "'0.9590"
$ws.Range("E5").Value = "  +1.37%  "

# Row 6
$ws.Range("D6").Value = "'280.53"
$ws.Range("E6").Value = "  +2.12%  "

# Row 7
$ws.Range("D7").Value = "'0.3614"
$ws.Range("E7").Value = "  -0.90%  "

# Row 8
$ws.Range("D8").Value = "'0.3165"
$ws.Range("E8").Value = "  +2.97%  "

# Row 9
$ws.Range("D9").Value = "'1.106"
$ws.Range("E9").Value = "  +7.00%  "

# Row 10
$ws.Range("D10").Value = "'40.46"
$ws.Range("E10").Value = "  +1.77%  "

# Row 11
$ws.Range("D11").Value = "'0.06764"
$ws.Range("E11").Value = "  +3.15%  "

# Row 12
$ws.Range("D12").Value = "'0.9978"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.639"
$ws.Range("E13").Value = "  +4.54%  "

# Row 14
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "'18.68"
$ws.Range("E14").Value = "  +3.93%  "

# Row 15
$ws.Range("D15").Value = "'6.306"
$ws.Range("E15").Value = "  +3.15%  "

# Row 16
$ws.Range("D16").Value = "'0.00001042"
$ws.Range("E16").Value = "  +1.83%  "

# Row 17
$ws.Range("D17").Value = "'0.9586"
$ws.Range("E17").Value = "  -0.40%  "

# Row 18
$ws.Range("D18").Value = "1.529.60"
$ws.Range("E18").Value = "  +4.94%  "

# Row 19
$ws.Range("D19").Value = "'0.06038"
$ws.Range("E19").Value = "  +4.95%  "

# Row 20
$ws.Range("D20").Value = "'71.52"
$ws.Range("E20").Value = "  +2.78%  "

# Row 21
$ws.Range("D21").Value = "'5.655"
$ws.Range("E21").Value = "  +4.42%  "

# Row 22
$ws.Range("D22").Value = "'14.95"
$ws.Range("E22").Value = "  +3.83%  "

# Row 23
$ws.Range("D23").Value = "'11.35"
$ws.Range("E23").Value = "  +4.57%  "

# Row 24
$ws.Range("E24").Value = "  +3.22%  "

# Row 25
$ws.Range("D25").Value = "21.054.29"
$ws.Range("E25").Value = "  +3.12%  "

# Row 26
$ws.Range("D26").Value = "'147.51"
$ws.Range("E26").Value = "  +4.61%  "

# Row 27
$ws.Range("D27").Value = "'2.196"
$ws.Range("E27").Value = "  +5.46%  "

# Row 28
$ws.Range("D28").Value = "'17.66"
$ws.Range("E28").Value = "  +3.29%  "

# Row 29
$ws.Range("D29").Value = "1.695.05"
$ws.Range("E29").Value = "  +5.24%  "

# Row 30
$ws.Range("D30").Value = "'117.86"
$ws.Range("E30").Value = "  +5.32%  "

# Row 31
$ws.Range("D31").Value = "'4.082"
$ws.Range("E31").Value = "  +7.11%  "

# Row 32
$ws.Range("D32").Value = "'0.8462"
$ws.Range("E32").Value = "  +7.58%  "

# Row 33
$ws.Range("D33").Value = "'5.152"
$ws.Range("E33").Value = "  +5.95%  "

# Row 34
$ws.Range("D34").Value = "'0.08015"
$ws.Range("E34").Value = "  +2.76%  "

# Row 35
$ws.Range("D35").Value = "'1.487"
$ws.Range("E35").Value = "  -1.04%  "

# Row 36
$ws.Range("D36").Value = "'1.212"
$ws.Range("E36").Value = "  +7.57%  "

# Row 37
$ws.Range("D37").Value = "'4.917"
$ws.Range("E37").Value = "  +5.69%  "

# Row 38
$ws.Range("D38").Value = "'0.05837"
$ws.Range("E38").Value = "  +2.57%  "

# Row 39
$ws.Range("D39").Value = "'0.02086"
$ws.Range("E39").Value = "  +2.81%  "

# Row 40
$ws.Range("D40").Value = "'10.64"
$ws.Range("E40").Value = "  +3.03%  "

# Row 41
$ws.Range("D41").Value = "'0.9589"
$ws.Range("E41").Value = "  +0.72%  "

# Row 42
$ws.Range("D42").Value = "'0.1900"
$ws.Range("E42").Value = "  +2.26%  "

# Row 43
$ws.Range("D43").Value = "'7.578"
$ws.Range("E43").Value = "  +2.28%  "

# Row 44
$ws.Range("D44").Value = "'0.5405"
$ws.Range("E44").Value = "  +2.85%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.39"
$ws.Range("E45").Value = "  +4.74%  "

# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.554"
$ws.Range("E46").Value = "  +1.99%  "

# Row 47
$ws.Range("D47").Value = "'0.5440"
$ws.Range("E47").Value = "  +5.91%  "

# Row 48
$ws.Range("D48").Value = "'121.27"
$ws.Range("E48").Value = "  +3.60%  "

# Row 49
$ws.Range("D49").Value = "'1.860"
$ws.Range("E49").Value = "  +6.56%  "

# Row 50
$ws.Range("D50").Value = "'0.06596"
$ws.Range("E50").Value = "  +2.67%  "

# Row 51
$ws.Range("D51").Value = "'0.9905"
$ws.Range("E51").Value = "  +0.31%  "
